$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (F1:H1) matching the existing header formatting
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$headerRng = $ws.Range("F1:H1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108   # xlCenter
$headerRng.VerticalAlignment = -4160     # xlTop
$headerRng.Borders.LineStyle = 1         # xlContinuous
$headerRng.Borders.Weight = 2            # xlThin

# New boolean columns: KNN_Outliers_MAD, SVM_Outliers_MAD, RF_Outliers_MAD
$ws.Range("F2").Value = $false
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false

$ws.Range("F3").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = $true

$ws.Range("F4").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = $false

$ws.Range("F5").Value = $false
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = $false
